# Applies the "make double changes again with main.docx file" edit:
#   1. A new, empty paragraph is inserted right after the first paragraph,
#      before "It will be treated as a binary file by Git."
#   2. The first paragraph's text gains two trailing spaces.
#   3. Three new dark-red (C00000) runs are appended to that same
#      first paragraph: "(This is a change – Version for branch ",
#      "main", ")".

$d = $word.ActiveDocument

# --- 1. Insert the new empty paragraph FIRST, while paragraph 1 still
#        only has its plain, uncolored text. Doing this before adding the
#        colored runs keeps the freshly inserted paragraph mark from
#        inheriting the dark-red color. ---
$d.Paragraphs(1).Range.InsertParagraphAfter()

# --- 2 & 3. Re-fetch paragraph 1 (still "This is a Microsoft word
#        document.") and grow it with the plain double space plus the
#        three colored runs. Each run gets its own InsertAfter call (on
#        the same, ever-growing Range object) so they stay as separate
#        <w:r> elements instead of merging with their neighbours. ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range

$r.InsertAfter("  ")

$darkRed = 192  # wdColor value for RGB C00000 (stored as BGR: R | G<<8 | B<<16)

# Paragraph Range.End always includes the trailing paragraph mark, so
# "$r.End - 1" is the real insertion point just before it.
$startA = $r.End - 1
$r.InsertAfter("(This is a change – Version for branch ")
$endA = $r.End - 1
$d.Range($startA, $endA).Font.Color = $darkRed

$startB = $r.End - 1
$r.InsertAfter("main")
$endB = $r.End - 1
$d.Range($startB, $endB).Font.Color = $darkRed

$startC = $r.End - 1
$r.InsertAfter(")")
$endC = $r.End - 1
$d.Range($startC, $endC).Font.Color = $darkRed
